$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7599915862083435
$ws.Range("B1").Value = 0.7077632546424866
$ws.Range("C1").Value = 3.883739709854126
$ws.Range("D1").Value = 3.049844026565552
$ws.Range("E1").Value = 0.8229544758796692
